$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555801496427"
$wb.Worksheets.Item(2).Name = "NB_TO-1651255582497778"
$wb.Worksheets.Item(3).Name = "RS_TO-1651255582497778"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555825457711"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555826161163"

# Sheet 1 (GNG) updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651255580114164.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555801326437.csv"
$ws1.Range("B4").Value = "go_stims-16512555801346438.csv"
$ws1.Range("B5").Value = "GNG_stims-1651255580148643.csv"

# Sheet 2 (NB) updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16512555815765684.csv"
$ws2.Range("B3").Value = "TB-16512555824821546.csv"
$ws2.Range("B4").Value = "ZB-match_5-16512555813586478.csv"
$ws2.Range("B5").Value = "OB-16512555818147686.csv"
$ws2.Range("B6").Value = "TB-16512555824665308.csv"
$ws2.Range("B7").Value = "ZB-match_5-16512555809341896.csv"
$ws2.Range("B8").Value = "TB-16512555820416002.csv"
$ws2.Range("B9").Value = "ZB-match_7-16512555802603693.csv"
$ws2.Range("B10").Value = "OB-16512555818872359.csv"

# Sheet 4 (TOL) updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555825134058.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255582497778.csv"
$ws4.Range("B4").Value = "MM_stims-16512555825290294.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555825134058.csv"
$ws4.Range("B6").Value = "MM_stims-16512555825447686.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555825290294.csv"

# Sheet 5 (vSAT) updates
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16512555825507696.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512555825848763.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555826004922.csv"
$ws5.Range("B5").Value = "SAT_stims-16512555825692403.csv"
